$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-04-13 Saturday" "2024-04-14 Sunday"

Replace-Text "426×2=852" "425×3=1275"
Replace-Text "386×9=3474" "187×7=1309"
Replace-Text "309×4=1236" "127×2=254"
Replace-Text "615×7=4305" "833×9=7497"
Replace-Text "525×5=2625" "671×3=2013"

Replace-Text "665×9=5985" "110×9=990"
Replace-Text "389×9=3501" "815×5=4075"
Replace-Text "345×8=2760" "507×6=3042"
Replace-Text "420×2=840" "788×5=3940"
Replace-Text "171×7=1197" "747×8=5976"

Replace-Text "716×7=5012" "754×7=5278"
Replace-Text "474×2=948" "641×9=5769"
Replace-Text "343×6=2058" "265×4=1060"
Replace-Text "614×2=1228" "627×3=1881"
Replace-Text "824×4=3296" "118×2=236"

Replace-Text "150×2=300" "510×9=4590"
Replace-Text "643×9=5787" "405×8=3240"
Replace-Text "849×2=1698" "690×4=2760"
Replace-Text "161×2=322" "739×3=2217"
Replace-Text "675×4=2700" "171×4=684"

Replace-Text "417×3=1251" "979×8=7832"
Replace-Text "543×5=2715" "504×9=4536"
Replace-Text "765×7=5355" "347×3=1041"
Replace-Text "447×7=3129" "797×3=2391"
Replace-Text "991×8=7928" "115×3=345"
